# Convention change to support multi-axle vehicles:
#   "sAxleF"     -> "sAxle1"
#   "Body_1Axle" -> "Body_Axle1"
# Applies to every worksheet in the workbook (Trailer_Elula,
# Trailer_Elula_Unstable, Trailer_Thwala) where column A row 5 holds the
# "sAxleF" label and column H row 4 holds the "Body_1Axle" class value.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("A5").Value2 -eq "sAxleF") {
        $ws.Range("A5").Value = "sAxle1"
    }
    if ($ws.Range("H4").Value2 -eq "Body_1Axle") {
        $ws.Range("H4").Value = "Body_Axle1"
    }
}
